$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 2171
$ws.Range("I58").Value = 198.33333
$ws.Range("J58").Value = 3247
$ws.Range("K58").Value = 594.99999
$ws.Range("L58").Value = 9741
$ws.Range("M58").Value = -444.99999
$ws.Range("N58").Value = -10041
# Row 125
$ws.Range("H125").Value = 1866.5
$ws.Range("I125").Value = 1432
$ws.Range("J125").Value = 1928.5714
$ws.Range("K125").Value = 12888
$ws.Range("L125").Value = 17357.1426
$ws.Range("M125").Value = -10428
$ws.Range("N125").Value = -22277.1426
# Row 129
$ws.Range("H129").Value = 1286.3846
$ws.Range("I129").Value = 547.75
$ws.Range("J129").Value = 1420.6818
$ws.Range("K129").Value = 1643.25
$ws.Range("L129").Value = 4262.0454
$ws.Range("M129").Value = 3356.75
# Row 131
$ws.Range("H131").Value = 3229.724
$ws.Range("I131").Value = 620.1111
$ws.Range("J131").Value = 7500
$ws.Range("K131").Value = 1860.3333
$ws.Range("L131").Value = 22500
$ws.Range("M131").Value = 3179.6667
$ws.Range("N131").Value = -32580
# Row 137
$ws.Range("H137").Value = 1357.0605
$ws.Range("I137").Value = 1330.7174
$ws.Range("J137").Value = 1417.65
$ws.Range("K137").Value = 3992.1522
$ws.Range("L137").Value = 4252.950000000001
$ws.Range("M137").Value = -1442.1522
$ws.Range("N137").Value = -9352.950000000001
# Row 138
$ws.Range("H138").Value = 1874.7142
$ws.Range("I138").Value = 777.75
$ws.Range("J138").Value = 4415.0527
$ws.Range("K138").Value = 2333.25
$ws.Range("L138").Value = 13245.1581
$ws.Range("M138").Value = 2806.75
$ws.Range("N138").Value = -23525.1581
# Row 140
$ws.Range("H140").Value = 79800
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 79800
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 79800
$ws.Range("N140").Value = -90160

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1656.481
$ws.Range("I74").Value = 1631.8823
$ws.Range("J74").Value = 1808.5454
$ws.Range("K74").Value = 1631.8823
$ws.Range("L74").Value = 1808.5454
$ws.Range("M74").Value = -757.8823
$ws.Range("N74").Value = -3556.5454
# Row 77
$ws.Range("H77").Value = 1656.481
$ws.Range("I77").Value = 1631.8823
$ws.Range("J77").Value = 1808.5454
$ws.Range("K77").Value = 8159.4115
$ws.Range("L77").Value = 9042.726999999999
$ws.Range("M77").Value = -3791.4115
$ws.Range("N77").Value = -17778.727
# Row 132
$ws.Range("H132").Value = 3717.7454
$ws.Range("I132").Value = 1545.6757
$ws.Range("J132").Value = 8182.5557
$ws.Range("K132").Value = 4637.0271
$ws.Range("L132").Value = 24547.6671
$ws.Range("M132").Value = -2107.0271
$ws.Range("N132").Value = -29607.6671

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 6
$ws.Range("H6").Value = 10333.333
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 10333.333
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 10333.333
$ws.Range("N6").Value = -10559.333
# Row 62
$ws.Range("H62").Value = 34250
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 34250
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 34250
$ws.Range("N62").Value = -35622
# Row 65
$ws.Range("H65").Value = 34250
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 34250
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 102750
$ws.Range("N65").Value = -109614

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 4468.125
$ws.Range("I7").Value = 7726.385
$ws.Range("J7").Value = 617.4545000000001
$ws.Range("K7").Value = 7726.385
$ws.Range("L7").Value = 617.4545000000001
$ws.Range("M7").Value = -7613.385
$ws.Range("N7").Value = -843.4545000000001
# Row 16
$ws.Range("H16").Value = 3150
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 3150
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 3150
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -3724
# Row 31
$ws.Range("H31").Value = 4763468.5
$ws.Range("I31").Value = 1130.6296
$ws.Range("J31").Value = 20836360
$ws.Range("K31").Value = 1130.6296
$ws.Range("L31").Value = 20836360
$ws.Range("M31").Value = -835.6296
$ws.Range("N31").Value = -20836950
# Row 34
$ws.Range("H34").Value = 4763468.5
$ws.Range("I34").Value = 1130.6296
$ws.Range("J34").Value = 20836360
$ws.Range("K34").Value = 1130.6296
$ws.Range("L34").Value = 20836360
$ws.Range("M34").Value = -928.6296
$ws.Range("N34").Value = -20836764
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 105
$ws.Range("H105").Value = 771204.7
$ws.Range("I105").Value = 835437.5
$ws.Range("J105").Value = 411
$ws.Range("K105").Value = 835437.5
$ws.Range("L105").Value = 411
$ws.Range("M105").Value = -833690.5
$ws.Range("N105").Value = -3905
# Row 113
$ws.Range("H113").Value = 3150
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3150
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3150
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -7490

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 725.4400000000001
$ws.Range("I5").Value = 513.8823
$ws.Range("J5").Value = 1175
$ws.Range("K5").Value = 1541.6469
$ws.Range("L5").Value = 3525
$ws.Range("M5").Value = -1429.6469
$ws.Range("N5").Value = -3749
# Row 18
$ws.Range("H18").Value = 1480.7646
$ws.Range("I18").Value = 521.8182
$ws.Range("J18").Value = 3238.8333
$ws.Range("K18").Value = 1565.4546
$ws.Range("L18").Value = 9716.499899999999
$ws.Range("M18").Value = -1396.4546
# Row 34
$ws.Range("H34").Value = 2028.65
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 2298.4119
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 6895.2357
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -7063.2357
# Row 80
$ws.Range("H80").Value = 2632.9092
$ws.Range("I80").Value = 994
$ws.Range("J80").Value = 3247.5
$ws.Range("K80").Value = 2982
$ws.Range("L80").Value = 9742.5
$ws.Range("M80").Value = -2046
$ws.Range("N80").Value = -11614.5
# Row 83
$ws.Range("H83").Value = 2632.9092
$ws.Range("I83").Value = 994
$ws.Range("J83").Value = 3247.5
$ws.Range("K83").Value = 8946
$ws.Range("L83").Value = 29227.5
$ws.Range("M83").Value = -4266
$ws.Range("N83").Value = -38587.5
# Row 86
$ws.Range("H86").Value = 1652
$ws.Range("I86").Value = 480
$ws.Range("J86").Value = 1945
$ws.Range("K86").Value = 1440
$ws.Range("L86").Value = 5835
$ws.Range("M86").Value = -254
# Row 89
$ws.Range("H89").Value = 1652
$ws.Range("I89").Value = 480
$ws.Range("J89").Value = 1945
$ws.Range("K89").Value = 4320
$ws.Range("L89").Value = 17505
$ws.Range("M89").Value = 1608
# Row 92
$ws.Range("H92").Value = 1364.5714
$ws.Range("I92").Value = 937
$ws.Range("J92").Value = 1685.25
$ws.Range("K92").Value = 2811
$ws.Range("L92").Value = 5055.75
$ws.Range("M92").Value = -1563
$ws.Range("N92").Value = -7551.75
# Row 109
$ws.Range("H109").Value = 4196.212
$ws.Range("I109").Value = 1162.5
$ws.Range("J109").Value = 4870.3706
$ws.Range("K109").Value = 3487.5
$ws.Range("L109").Value = 14611.1118
$ws.Range("M109").Value = -2447.5
$ws.Range("N109").Value = -16691.1118
# Row 113
$ws.Range("H113").Value = 505.22223
$ws.Range("I113").Value = 492.22726
$ws.Range("J113").Value = 525.6429000000001
$ws.Range("K113").Value = 1476.68178
$ws.Range("L113").Value = 1576.9287
$ws.Range("M113").Value = 693.3182200000001
$ws.Range("N113").Value = -5916.9287
# Row 122
$ws.Range("H122").Value = 2277.6978
$ws.Range("I122").Value = 253.91667
$ws.Range("J122").Value = 3061.0967
$ws.Range("K122").Value = 2285.25003
$ws.Range("L122").Value = 27549.8703
$ws.Range("M122").Value = 164.7499699999998
$ws.Range("N122").Value = -32449.8703
# Row 129
$ws.Range("H129").Value = 2082.4138
$ws.Range("I129").Value = 554.1667
$ws.Range("J129").Value = 3161.1765
$ws.Range("K129").Value = 1662.5001
$ws.Range("L129").Value = 9483.529500000001
$ws.Range("M129").Value = 3337.4999
$ws.Range("N129").Value = -19483.5295
# Row 132
$ws.Range("H132").Value = 2317
$ws.Range("I132").Value = 971.2857
$ws.Range("J132").Value = 4201
$ws.Range("K132").Value = 8741.5713
$ws.Range("L132").Value = 37809
$ws.Range("M132").Value = -6211.5713
$ws.Range("N132").Value = -42869
# Row 135
$ws.Range("H135").Value = 725.4400000000001
$ws.Range("I135").Value = 513.8823
$ws.Range("J135").Value = 1175
$ws.Range("K135").Value = 4624.9407
$ws.Range("L135").Value = 10575
$ws.Range("M135").Value = -2089.9407
$ws.Range("N135").Value = -15645

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -5826
# Row 42
$ws.Range("H42").Value = 70049
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 70049
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 70049
$ws.Range("N42").Value = -70805
# Row 126
$ws.Range("H126").Value = 1980.5526
$ws.Range("I126").Value = 2447.4285
$ws.Range("J126").Value = 1403.8235
$ws.Range("K126").Value = 7342.2855
$ws.Range("L126").Value = 4211.470499999999
$ws.Range("M126").Value = -4872.2855
$ws.Range("N126").Value = -9151.470499999999
